$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.9146506000311661
$ws.Range("C2").Value = 0.2517570627570933
$ws.Range("D2").Value = 0.2206719225841596
$ws.Range("F2").Value = 1.510899119732557
$ws.Range("G2").Value = 0.8619426395293885
$ws.Range("H2").Value = 0.8994347571257322
$ws.Range("I2").Value = 0.6558943847060199
$ws.Range("J2").Value = 0.2570088181630794
$ws.Range("M2").Value = 0.449619215129033

# Row 3
$ws.Range("B3").Value = 0.8163795911088414
$ws.Range("C3").Value = 0.2220404484730238
$ws.Range("D3").Value = 0.2188549481036191
$ws.Range("F3").Value = 1.509394915658518
$ws.Range("G3").Value = 0.8580574676599326
$ws.Range("H3").Value = 0.9037903686295152
$ws.Range("I3").Value = 0.6634440360821685
$ws.Range("J3").Value = 0.2577166705739913
$ws.Range("M3").Value = 0.4216201576306986

# Row 4
$ws.Range("B4").Value = 0.7560858235492844
$ws.Range("C4").Value = 0.2037933306949355
$ws.Range("D4").Value = 0.2177977916562455
$ws.Range("F4").Value = 1.509570063004972
$ws.Range("G4").Value = 0.8564663236758889
$ws.Range("H4").Value = 0.9070495324632333
$ws.Range("I4").Value = 0.6686095207763891
$ws.Range("J4").Value = 0.2582897348443254
$ws.Range("M4").Value = 0.4045764733676691

# Row 5
$ws.Range("B5").Value = 0.7315278501415605
$ws.Range("C5").Value = 0.196357367823623
$ws.Range("D5").Value = 0.2173817571725962
$ws.Range("F5").Value = 1.509917058924188
$ws.Range("G5").Value = 0.8560168332202949
$ws.Range("H5").Value = 0.9085244621660422
$ws.Range("I5").Value = 0.6708473717683106
$ws.Range("J5").Value = 0.2585580392976539
$ws.Range("M5").Value = 0.3976684166713866

# Row 6
$ws.Range("B6").Value = 0.7274507828384742
$ws.Range("C6").Value = 0.1951226291052421
$ws.Range("D6").Value = 0.2173135685864409
$ws.Range("F6").Value = 1.509991304826158
$ws.Range("G6").Value = 0.8559541850861052
$ws.Range("H6").Value = 0.9087782315531143
$ws.Range("I6").Value = 0.6712269790398828
$ws.Range("J6").Value = 0.2586046904018815
$ws.Range("M6").Value = 0.3965236031235193

# Row 7
$ws.Range("B7").Value = 0.7557545754155512
$ws.Range("C7").Value = 0.2036930469755589
$ws.Range("D7").Value = 0.2177921210025175
$ws.Range("F7").Value = 1.509573627616859
$ws.Range("G7").Value = 0.8564594574094997
$ws.Range("H7").Value = 0.9070688298668017
$ws.Range("I7").Value = 0.6686391637381277
$ws.Range("J7").Value = 0.2582932125320667
$ws.Range("M7").Value = 0.4044831572606498

# Row 8
$ws.Range("B8").Value = 0.8807579162989327
$ws.Range("C8").Value = 0.2415110521550332
$ws.Range("D8").Value = 0.2200333296369763
$ws.Range("F8").Value = 1.510152022913005
$ws.Range("G8").Value = 0.8604376541126584
$ws.Range("H8").Value = 0.900815039204943
$ws.Range("I8").Value = 0.6583872745862998
$ws.Range("J8").Value = 0.2572241303712346
$ws.Range("M8").Value = 0.4399345658039948

# Row 9
$ws.Range("B9").Value = 1.12621982919768
$ws.Range("C9").Value = 0.3156621019641932
$ws.Range("D9").Value = 0.2248899687876929
$ws.Range("F9").Value = 1.520036616587689
$ws.Range("G9").Value = 0.8745813624786791
$ws.Range("H9").Value = 0.8932039334132043
$ws.Range("I9").Value = 0.642506106533304
$ws.Range("J9").Value = 0.2562279817626347
$ws.Range("M9").Value = 0.510622880619735

# Row 10
$ws.Range("B10").Value = 1.306746066568905
$ws.Range("C10").Value = 0.3701379283293136
$ws.Range("D10").Value = 0.2287370222664009
$ws.Range("F10").Value = 1.532681674975066
$ws.Range("G10").Value = 0.8888966043161162
$ws.Range("H10").Value = 0.8904666673672637
$ws.Range("I10").Value = 0.6334371152161395
$ws.Range("J10").Value = 0.2561698249073316
$ws.Range("M10").Value = 0.5632691966780214

# Row 11
$ws.Range("B11").Value = 1.388910436049741
$ws.Range("C11").Value = 0.3949208618375337
$ws.Range("D11").Value = 0.2305472117321727
$ws.Range("F11").Value = 1.539613746467296
$ws.Range("G11").Value = 0.8962735084505624
$ws.Range("H11").Value = 0.8898453000640671
$ws.Range("I11").Value = 0.6298810390707104
$ws.Range("J11").Value = 0.2562903211923953
$ws.Range("M11").Value = 0.5873742354146145

# Row 12
$ws.Range("B12").Value = 1.42002943852674
$ws.Range("C12").Value = 0.404305716606757
$ws.Range("D12").Value = 0.2312412811324833
$ws.Range("F12").Value = 1.542409193915546
$ws.Range("G12").Value = 0.8991922734651894
$ws.Range("H12").Value = 0.8897000241273645
$ws.Range("I12").Value = 0.6286167682313462
$ws.Range("J12").Value = 0.2563571305825789
$ws.Range("M12").Value = 0.5965245201947198

# Row 13
$ws.Range("B13").Value = 1.41332719030612
$ws.Range("C13").Value = 0.4022845180568311
$ws.Range("D13").Value = 0.2310914197587692
$ws.Range("F13").Value = 1.54179955237575
$ws.Range("G13").Value = 0.8985580779966824
$ws.Range("H13").Value = 0.8897273030724051
$ws.Range("I13").Value = 0.6288853823767369
$ws.Range("J13").Value = 0.2563417992183332
$ws.Range("M13").Value = 0.5945528570658638

# Row 14
$ws.Range("B14").Value = 1.391470517273262
$ws.Range("C14").Value = 0.395692958324787
$ws.Range("D14").Value = 0.2306041414185529
$ws.Range("F14").Value = 1.539840309657862
$ws.Range("G14").Value = 0.8965111206444334
$ws.Range("H14").Value = 0.8898315421252789
$ws.Range("I14").Value = 0.6297753739840388
$ws.Range("J14").Value = 0.2562953928818814
$ws.Range("M14").Value = 0.5881265904469473

# Row 15
$ws.Range("B15").Value = 1.378083316458003
$ws.Range("C15").Value = 0.3916554462016961
$ws.Range("D15").Value = 0.2303067862932409
$ws.Range("F15").Value = 1.538662433924387
$ws.Range("G15").Value = 0.8952736444413318
$ws.Range("H15").Value = 0.8899071245435977
$ws.Range("I15").Value = 0.6303312557875245
$ws.Range("J15").Value = 0.2562697273783385
$ws.Range("M15").Value = 0.5841932057795418

# Row 16
$ws.Range("B16").Value = 1.301377200118452
$ws.Range("C16").Value = 0.3685183259796077
$ws.Range("D16").Value = 0.2286199273558225
$ws.Range("F16").Value = 1.532252453817392
$ws.Range("G16").Value = 0.8884319883389225
$ws.Range("H16").Value = 0.890519855239404
$ws.Range("I16").Value = 0.6336810105782362
$ws.Range("J16").Value = 0.2561649113810134
$ws.Range("M16").Value = 0.5616969970901806

# Row 17
$ws.Range("B17").Value = 1.254330604519168
$ws.Range("C17").Value = 0.354324765657509
$ws.Range("D17").Value = 0.2276004565765106
$ws.Range("F17").Value = 1.528622855461492
$ws.Range("G17").Value = 0.884457018849389
$ws.Range("H17").Value = 0.8910557479924677
$ws.Range("I17").Value = 0.6358821482501682
$ws.Range("J17").Value = 0.2561382825022136
$ws.Range("M17").Value = 0.5479361056786161

# Row 18
$ws.Range("B18").Value = 1.227274660974672
$ws.Range("C18").Value = 0.3461611641700983
$ws.Range("D18").Value = 0.2270197486065939
$ws.Range("F18").Value = 1.526646205763285
$ws.Range("G18").Value = 0.8822520676351644
$ws.Range("H18").Value = 0.8914226780153882
$ws.Range("I18").Value = 0.6372017570448989
$ws.Range("J18").Value = 0.2561367951751095
$ws.Range("M18").Value = 0.5400358966005285

# Row 19
$ws.Range("B19").Value = 1.218114700875674
$ws.Range("C19").Value = 0.343397140495199
$ws.Range("D19").Value = 0.2268241055830771
$ws.Range("F19").Value = 1.525995986870981
$ws.Range("G19").Value = 0.881519451992105
$ws.Range("H19").Value = 0.8915569861562886
$ws.Range("I19").Value = 0.6376577413244178
$ws.Range("J19").Value = 0.2561386652475761
$ws.Range("M19").Value = 0.5373635509828887

# Row 20
$ws.Range("B20").Value = 1.259338386181867
$ws.Range("C20").Value = 0.3558356789285426
$ws.Range("D20").Value = 0.2277083951408372
$ws.Range("F20").Value = 1.528997738796377
$ws.Range("G20").Value = 0.8848717353828164
$ws.Range("H20").Value = 0.8909926239057739
$ws.Range("I20").Value = 0.6356422855501904
$ws.Range("J20").Value = 0.2561396856025908
$ws.Range("M20").Value = 0.5493994569777954

# Row 21
$ws.Range("B21").Value = 1.39789021965646
$ws.Range("C21").Value = 0.3976290576536599
$ws.Range("D21").Value = 0.2307470342157103
$ws.Range("F21").Value = 1.540411155350753
$ws.Range("G21").Value = 0.8971089538428743
$ws.Range("H21").Value = 0.8897984788768269
$ws.Range("I21").Value = 0.6295117236446259
$ws.Range("J21").Value = 0.2563084483377196
$ws.Range("M21").Value = 0.5900135389687335

# Row 22
$ws.Range("B22").Value = 1.488471396760531
$ws.Range("C22").Value = 0.4249439455597326
$ws.Range("D22").Value = 0.2327830024510007
$ws.Range("F22").Value = 1.548864195462428
$ws.Range("G22").Value = 0.9058374109048373
$ws.Range("H22").Value = 0.8895429033086657
$ws.Range("I22").Value = 0.6259851859963277
$ws.Range("J22").Value = 0.2565422210220376
$ws.Range("M22").Value = 0.6166866734702268

# Row 23
$ws.Range("B23").Value = 1.440124110674276
$ws.Range("C23").Value = 0.4103654734738598
$ws.Range("D23").Value = 0.2316918083075734
$ws.Range("F23").Value = 1.544261468163228
$ws.Range("G23").Value = 0.9011116998513131
$ws.Range("H23").Value = 0.889631178451765
$ws.Range("I23").Value = 0.6278232834672934
$ws.Range("J23").Value = 0.2564061378738174
$ws.Range("M23").Value = 0.6024389357550888

# Row 24
$ws.Range("B24").Value = 1.25707439366181
$ws.Range("C24").Value = 0.3551526060452375
$ws.Range("D24").Value = 0.2276595792944534
$ws.Range("F24").Value = 1.52882791127486
$ws.Range("G24").Value = 0.884683991937095
$ws.Range("H24").Value = 0.8910209790420396
$ws.Range("I24").Value = 0.6357505588435615
$ws.Range("J24").Value = 0.2561390082102619
$ws.Range("M24").Value = 0.548737841234896

# Row 25
$ws.Range("B25").Value = 1.059782226245147
$ws.Range("C25").Value = 0.2956030416919475
$ws.Range("D25").Value = 0.2235269639924411
$ws.Range("F25").Value = 1.516420790851441
$ws.Range("G25").Value = 0.8700701873592607
$ws.Range("H25").Value = 0.894762965086727
$ws.Range("I25").Value = 0.6463477677836593
$ws.Range("J25").Value = 0.2563793685274973
$ws.Range("M25").Value = 0.491374832376934

